$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stack")

# Existing hyperlinks (C2/C3) point at cells that are about to shift down one
# row; remove them now and re-add them against their final cell addresses
# once the new row is in place, instead of relying on the engine to shift
# the hyperlink anchors for us.
$ws.Hyperlinks.Delete()

# Insert a new row above the current row 2 and pull the formatting for the
# new row down from what is now row 3 (same "shape" as every other data row).
$ws.Rows.Item(2).Insert()
$ws.Range("A3:K3").Copy()
$ws.Range("A2:K2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new "Maximum Frequency Stack" entry.
$ws.Range("A2").Value = 895
$ws.Range("B2").Value = "Maximum Frequency Stack"
$ws.Range("C2").Value = "https://leetcode.com/problems/maximum-frequency-stack/description/"
$ws.Range("D2").Value = "Hard"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "Binary Search, Hash"
$ws.Range("G2").Value = "O(1)"
$ws.Range("H2").Value = "O(n)"
$ws.Range("I2").Value = 45512
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""

# Re-create the three hyperlinks at their final positions.
$ws.Hyperlinks.Add($ws.Range("C2"), "https://leetcode.com/problems/maximum-frequency-stack/description/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://leetcode.com/problems/dinner-plate-stacks/description/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://leetcode.com/problems/next-greater-element-iv/description/") | Out-Null

# Hyperlinks.Add resets the run font to the generic "Hyperlink" style; put the
# sheet's local hyperlink font back so the three links look the same.
$ws.Range("C2:C4").Font.Name = "Microsoft YaHei"
$ws.Range("C2:C4").Font.Size = 10

# Extend the Easy/Medium/Hard conditional formatting down to the new row
# without losing the dxf formatting already tied to these rules.
$fcs = $ws.Range("D2:D3").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
  $fcs.Item($i).ModifyAppliesToRange($ws.Range("D2:D4"))
}

Write-Output "done"
